$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data rows
$ws.Range("A3").Value = "/404"
$ws.Range("B3").Value = "/nextgen"
$ws.Range("A4").Value = "/contactus"
$ws.Range("B4").Value = "/home"

# Apply gray fill ("White, Background 1, Darker 15%") to header row A1:C1
$headerRange = $ws.Range("A1:C1")
$headerRange.Interior.Color = 14277081

# Update selection to C9
$ws.Range("C9").Select()
